$wb = $excel.ActiveWorkbook

# Per-sheet test-run data: Result (col A) + Date (col B) for data rows 2-8.
# Order below matches the order the runs completed (so new shared strings land
# in the same sequence as the authored workbook).
$sheetData = @(
    @{ Name = "Debit"; Result = "Pass"; Dates = @("Tue Feb 28 18:36:43 EST 2023", "Tue Feb 28 18:36:58 EST 2023", "Tue Feb 28 18:37:12 EST 2023", "Tue Feb 28 18:37:27 EST 2023", "Tue Feb 28 18:37:41 EST 2023", "Tue Feb 28 18:37:55 EST 2023", "Tue Feb 28 18:38:10 EST 2023") }
    @{ Name = "Debit-ZeroDollar"; Result = "Pass"; Dates = @("Tue Feb 28 18:38:24 EST 2023", "Tue Feb 28 18:38:39 EST 2023", "Tue Feb 28 18:38:53 EST 2023", "Tue Feb 28 18:39:07 EST 2023", "Tue Feb 28 18:39:21 EST 2023", "Tue Feb 28 18:39:35 EST 2023", "Tue Feb 28 18:39:50 EST 2023") }
    @{ Name = "Debit-Void"; Result = "Pass"; Dates = @("Tue Feb 28 18:40:04 EST 2023", "Tue Feb 28 18:40:28 EST 2023", "Tue Feb 28 18:40:52 EST 2023", "Tue Feb 28 18:41:16 EST 2023", "Tue Feb 28 18:41:40 EST 2023", "Tue Feb 28 18:42:04 EST 2023", "Tue Feb 28 18:42:27 EST 2023") }
    @{ Name = "Debit-Credit"; Result = "Pass"; Dates = @("Tue Feb 28 18:42:51 EST 2023", "Tue Feb 28 18:43:16 EST 2023", "Tue Feb 28 18:43:41 EST 2023", "Tue Feb 28 18:44:06 EST 2023", "Tue Feb 28 18:44:31 EST 2023", "Tue Feb 28 18:44:56 EST 2023", "Tue Feb 28 18:45:21 EST 2023") }
    @{ Name = "Debit-Credit-Void"; Result = "Pass"; Dates = @("Tue Feb 28 18:45:46 EST 2023", "Tue Feb 28 18:46:20 EST 2023", "Tue Feb 28 18:46:55 EST 2023", "Tue Feb 28 18:47:30 EST 2023", "Tue Feb 28 18:48:03 EST 2023", "Tue Feb 28 18:48:37 EST 2023", "Tue Feb 28 18:49:12 EST 2023") }
    @{ Name = "DebitCredit-RemID-Pipe"; Result = "Pass"; Dates = @("Tue Feb 28 18:49:48 EST 2023", "Tue Feb 28 18:50:13 EST 2023", "Tue Feb 28 18:50:37 EST 2023", "Tue Feb 28 18:51:02 EST 2023", "Tue Feb 28 18:51:27 EST 2023", "Tue Feb 28 18:51:52 EST 2023", "Tue Feb 28 18:52:17 EST 2023") }
    @{ Name = "Debit-RemID-Pipe"; Result = "Pass"; Dates = @("Tue Feb 28 18:52:42 EST 2023", "Tue Feb 28 18:52:57 EST 2023", "Tue Feb 28 18:53:11 EST 2023", "Tue Feb 28 18:53:25 EST 2023", "Tue Feb 28 18:53:41 EST 2023", "Tue Feb 28 18:53:56 EST 2023", "Tue Feb 28 18:54:11 EST 2023") }
    @{ Name = "DebitVoid-RemID-Pipe"; Result = "Pass"; Dates = @("Tue Feb 28 18:54:26 EST 2023", "Tue Feb 28 18:54:50 EST 2023", "Tue Feb 28 18:55:14 EST 2023", "Tue Feb 28 18:55:38 EST 2023", "Tue Feb 28 18:56:01 EST 2023", "Tue Feb 28 18:56:25 EST 2023", "Tue Feb 28 18:56:49 EST 2023") }
)

foreach ($entry in $sheetData) {
    $ws = $wb.Worksheets.Item($entry.Name)
    for ($i = 0; $i -lt $entry.Dates.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $entry.Result
        $ws.Cells.Item($row, 2).Value = $entry.Dates[$i]
    }
}
